$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell C1: "hi" -> "thanks"
$ws.Cells.Item(1, 3).Value = "thanks"

# A13: "мис   спасибо" -> "спасибо"
$ws.Cells.Item(13, 1).Value = "спасибо"
